$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale alert rows (3-7); only the corrected alert (row 2) remains
$ws.Range("A3:J7").EntireRow.Delete()

# Update the remaining alert row with the corrected timestamp and raw JSON payload
$ws.Range("A2").Value = "2025-11-25T19:02:41.814Z"
$ws.Range("J2").Value = '{"Sequence":1,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"a3f6808f-a153-4ec0-901d-867e2662dacd","EventDtm":"2025-11-25T19:02:25Z","AppDtm":"2025-11-25T19:02:41Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T19:02:25Z","DeviceID":"JSGA621320217","GPSLockState":"LOCKED","GPSSatelliteCount":20,"GPSLastLock":0,"GPSLatitude":26.644599,"GPSLongitude":88.467988,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"405","MNC":"51","LAC":"41512","CID":"16323942","RSSI":null,"ExtPower":true,"ExtPowerVoltage":34,"BatteryVoltage":8.1,"DeviceTemp":25,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CXRU1036301","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-25T19:02:25Z","AssetType":"Reefer","AssetID":"CXRU1036301","OEM":"DAIKIN","TAmb":18.81,"TAmbQ":null,"TUSDA4":-53.4,"TUSDA4Q":"OOR","PctCO2":25.5,"PctCO2Q":"OOR","PctCO2Set":25.5,"PctCO2SetQ":"OOR","PSuc":-10,"PSucQ":"asProvided","TDis":64.3,"TDisQ":null,"FreqComp":null,"TSuc":-24.3,"TSucQ":null,"MCond":"On","PCond":null,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":null,"SnCtrl":null,"AmpPhA":null,"AmpPhB":null,"AmpPhC":null,"PDis":830,"PDisQ":"asProvided","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":25,"PctExpVlvQ":null,"TEvap":-1.38,"TEvapQ":null,"MCtrl3":"Modulation","PctHtr":null,"PctHtrQ":null,"MEvapFanHS":"On","PctGasVlv":0,"PctGasVlvQ":null,"PctHum":100.39,"PctHumQ":"OOR","PctHumSet":65,"PctHumSetQ":"off","FreqLine":null,"FreqLineQ":null,"VLine1":424.8,"VLine2":null,"VLine3":null,"MEvapFanLS":"Off","PctO2":25.5,"PctO2Q":"OOR","PctO2Set":25.5,"PctO2SetQ":"OOR","MCtrl2":"Modulation","TRtn1":1.19,"TRtn1Q":null,"TRtn2":1.9,"TRtn2Q":null,"TSet":0,"TSetQ":null,"VerSWMajor":"265A","VerSWMinor":null,"PctSucVlv":12.8,"PctSucVlvQ":null,"TSup1":0.19,"TSup1Q":null,"TSup2":1.8,"TSup2Q":null,"AmpTotalDraw":13,"AmpTotalDrawQ":"asProvided","TUSDA1":-53.4,"TUSDA1Q":"OOR","TUSDA2":-53.4,"TUSDA2Q":"OOR","TUSDA3":-53.4,"TUSDA3Q":"OOR","CmhVent":1020,"CmhVentQ":"error","BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-07-16T10:29:07Z","PTIResult":"Failed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":"UnitOff","ReeferAlarms":null}}}'

Write-Host "Alerts sheet now has" $ws.UsedRange.Rows.Count "row(s) of data (incl. header)"
